{"js": "// Replace the division-problem text in each table cell according to the\n// fixed old->new mapping below. The mapping is applied positionally, in\n// document order, because some of the \"before\" values are not unique once\n// combined with the \"after\" values (e.g. two cells both become \"41\u00f74=\").\nconst replacements = [\n  [\"68\u00f79=\", \"89\u00f74=\"],\n  [\"43\u00f75=\", \"76\u00f75=\"],\n  [\"97\u00f79=\", \"97\u00f77=\"],\n  [\"55\u00f77=\", \"29\u00f75=\"],\n  [\"23\u00f72=\", \"65\u00f74=\"],\n  [\"55\u00f76=\", \"67\u00f79=\"],\n  [\"78\u00f78=\", \"18\u00f73=\"],\n  [\"58\u00f72=\", \"17\u00f77=\"],\n  [\"90\u00f73=\", \"41\u00f74=\"],\n  [\"93\u00f76=\", \"41\u00f74=\"],\n  [\"54\u00f77=\", \"96\u00f74=\"],\n  [\"38\u00f77=\", \"14\u00f76=\"],\n  [\"52\u00f78=\", \"17\u00f79=\"],\n  [\"94\u00f77=\", \"87\u00f72=\"],\n  [\"52\u00f74=\", \"32\u00f74=\"],\n  [\"89\u00f74=\", \"65\u00f75=\"],\n  [\"26\u00f73=\", \"86\u00f76=\"],\n  [\"97\u00f75=\", \"71\u00f74=\"],\n  [\"10\u00f73=\", \"41\u00f77=\"],\n  [\"44\u00f78=\", \"48\u00f74=\"],\n  [\"62\u00f79=\", \"91\u00f78=\"],\n  [\"58\u00f79=\", \"62\u00f75=\"],\n  [\"68\u00f78=\", \"63\u00f75=\"],\n  [\"57\u00f79=\", \"46\u00f75=\"],\n  [\"20\u00f72=\", \"95\u00f74=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet replIndex = 0;\nfor (let i = 0; i < paragraphs.items.length && replIndex < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const current = para.text;\n  const [oldText, newText] = replacements[replIndex];\n  if (current === oldText) {\n    para.insertText(newText, \"Replace\");\n    replIndex++;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Fixed old->new mapping, applied positionally in document order (some\n# \"before\" values are not globally unique once combined with \"after\"\n# values, e.g. two cells both become \"41\u00f74=\", and \"89\u00f74=\" is both a\n# pre-existing value and the result of another cell's replacement).\n$replacements = @(\n    ,@(\"68\u00f79=\", \"89\u00f74=\")\n    ,@(\"43\u00f75=\", \"76\u00f75=\")\n    ,@(\"97\u00f79=\", \"97\u00f77=\")\n    ,@(\"55\u00f77=\", \"29\u00f75=\")\n    ,@(\"23\u00f72=\", \"65\u00f74=\")\n    ,@(\"55\u00f76=\", \"67\u00f79=\")\n    ,@(\"78\u00f78=\", \"18\u00f73=\")\n    ,@(\"58\u00f72=\", \"17\u00f77=\")\n    ,@(\"90\u00f73=\", \"41\u00f74=\")\n    ,@(\"93\u00f76=\", \"41\u00f74=\")\n    ,@(\"54\u00f77=\", \"96\u00f74=\")\n    ,@(\"38\u00f77=\", \"14\u00f76=\")\n    ,@(\"52\u00f78=\", \"17\u00f79=\")\n    ,@(\"94\u00f77=\", \"87\u00f72=\")\n    ,@(\"52\u00f74=\", \"32\u00f74=\")\n    ,@(\"89\u00f74=\", \"65\u00f75=\")\n    ,@(\"26\u00f73=\", \"86\u00f76=\")\n    ,@(\"97\u00f75=\", \"71\u00f74=\")\n    ,@(\"10\u00f73=\", \"41\u00f77=\")\n    ,@(\"44\u00f78=\", \"48\u00f74=\")\n    ,@(\"62\u00f79=\", \"91\u00f78=\")\n    ,@(\"58\u00f79=\", \"62\u00f75=\")\n    ,@(\"68\u00f78=\", \"63\u00f75=\")\n    ,@(\"57\u00f79=\", \"46\u00f75=\")\n    ,@(\"20\u00f72=\", \"95\u00f74=\")\n)\n\n$replIndex = 0\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count -and $replIndex -lt $replacements.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    # Paragraph.Range.Text carries a trailing paragraph mark (CR) and, for a\n    # paragraph that is also the last one in a table cell, a trailing cell\n    # mark (BEL) as well - strip both before comparing/replacing so only the\n    # visible text is touched.\n    $current = $r.Text.TrimEnd([char]13, [char]7)\n    $pair = $replacements[$replIndex]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    if ($current -eq $oldText) {\n        $r.Text = $newText\n        $replIndex++\n    }\n}\n\n$replIndex\n"}
